# Auto-generated edit script: updates cached market-price / profit values
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets (Leve profit tracker).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC (108 cell updates) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 2100.2917
$ws.Range("I33").Value = 2382
$ws.Range("K33").Value = 2382
$ws.Range("M33").Value = -2153
$ws.Range("H40").Value = 4498.1665
$ws.Range("J40").Value = 5499.5
$ws.Range("L40").Value = 5499.5
$ws.Range("N40").Value = -5849.5
$ws.Range("H69").Value = 13772.111
$ws.Range("J69").Value = 13472.5
$ws.Range("L69").Value = 40417.5
$ws.Range("N69").Value = -42165.5
$ws.Range("H72").Value = 13772.111
$ws.Range("J72").Value = 13472.5
$ws.Range("L72").Value = 121252.5
$ws.Range("N72").Value = -129988.5
$ws.Range("H74").Value = 18590.625
$ws.Range("I74").Value = 7908.3335
$ws.Range("K74").Value = 7908.3335
$ws.Range("M74").Value = -6972.3335
$ws.Range("H77").Value = 18590.625
$ws.Range("I77").Value = 7908.3335
$ws.Range("K77").Value = 39541.6675
$ws.Range("M77").Value = -34861.6675
$ws.Range("H80").Value = 298.30768
$ws.Range("I80").Value = 174.22223
$ws.Range("J80").Value = 577.5
$ws.Range("K80").Value = 522.66669
$ws.Range("L80").Value = 1732.5
$ws.Range("M80").Value = 475.33331
$ws.Range("N80").Value = -3728.5
$ws.Range("H83").Value = 298.30768
$ws.Range("I83").Value = 174.22223
$ws.Range("J83").Value = 577.5
$ws.Range("K83").Value = 1568.00007
$ws.Range("L83").Value = 5197.5
$ws.Range("M83").Value = 3423.99993
$ws.Range("N83").Value = -15181.5
$ws.Range("H92").Value = 30380.9
$ws.Range("I92").Value = 168.45
$ws.Range("J92").Value = 90805.8
$ws.Range("K92").Value = 168.45
$ws.Range("L92").Value = 90805.8
$ws.Range("M92").Value = 1079.55
$ws.Range("N92").Value = -93301.8
$ws.Range("H94").Value = 1952.5
$ws.Range("I94").Value = 1952.5
$ws.Range("K94").Value = 1952.5
$ws.Range("M94").Value = -1501.5
$ws.Range("H98").Value = 135083.28
$ws.Range("J98").Value = 6899
$ws.Range("L98").Value = 6899
$ws.Range("N98").Value = -9895
$ws.Range("H113").Value = 1544.7778
$ws.Range("I113").Value = 1657.5714
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 1657.5714
$ws.Range("L113").Value = 1150
$ws.Range("M113").Value = 1596.4286
$ws.Range("N113").Value = -7658
$ws.Range("H116").Value = 6401.185
$ws.Range("I116").Value = 5665.1875
$ws.Range("K116").Value = 5665.1875
$ws.Range("M116").Value = -2223.1875
$ws.Range("H122").Value = 135083.28
$ws.Range("J122").Value = 6899
$ws.Range("L122").Value = 20697
$ws.Range("N122").Value = -25597
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H129").Value = 6804
$ws.Range("I129").Value = 3402.889
$ws.Range("K129").Value = 10208.667
$ws.Range("M129").Value = -5208.667000000001
$ws.Range("H131").Value = 9816.280000000001
$ws.Range("I131").Value = 8876.5
$ws.Range("J131").Value = 11487
$ws.Range("K131").Value = 26629.5
$ws.Range("L131").Value = 34461
$ws.Range("M131").Value = -21589.5
$ws.Range("N131").Value = -44541
$ws.Range("H132").Value = 2535.8572
$ws.Range("I132").Value = 2535.8572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7607.571599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5077.571599999999
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 192454.03
$ws.Range("I137").Value = 1078.6061
$ws.Range("J137").Value = 3350148.5
$ws.Range("K137").Value = 3235.8183
$ws.Range("L137").Value = 10050445.5
$ws.Range("M137").Value = -685.8182999999999
$ws.Range("N137").Value = -10055545.5
$ws.Range("H138").Value = 2574.1086
$ws.Range("I138").Value = 1725.069
$ws.Range("K138").Value = 5175.207
$ws.Range("M138").Value = -35.20700000000033
$ws.Range("H141").Value = 5971.273
$ws.Range("I141").Value = 4503.3423
$ws.Range("J141").Value = 15268.167
$ws.Range("K141").Value = 13510.0269
$ws.Range("L141").Value = 45804.501
$ws.Range("M141").Value = -8330.026900000001
$ws.Range("N141").Value = -56164.501

# --- Sheet 2: ARM (101 cell updates) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 2016.4231
$ws.Range("I2").Value = 2004.3334
$ws.Range("J2").Value = 2043.625
$ws.Range("K2").Value = 2004.3334
$ws.Range("L2").Value = 2043.625
$ws.Range("M2").Value = -1891.3334
$ws.Range("N2").Value = -2269.625
$ws.Range("H4").Value = 2112.4
$ws.Range("I4").Value = 1983
$ws.Range("J4").Value = 2445.1428
$ws.Range("K4").Value = 1983
$ws.Range("L4").Value = 2445.1428
$ws.Range("M4").Value = -1867
$ws.Range("N4").Value = -2677.1428
$ws.Range("H32").Value = 3993.2856
$ws.Range("I32").Value = 3716.4238
$ws.Range("K32").Value = 3716.4238
$ws.Range("M32").Value = -3429.4238
$ws.Range("H38").Value = 3009.5
$ws.Range("I38").Value = 1019
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 1019
$ws.Range("L38").Value = 5000
$ws.Range("M38").Value = -552
$ws.Range("N38").Value = -5934
$ws.Range("H45").Value = 1878.579
$ws.Range("I45").Value = 1356
$ws.Range("K45").Value = 1356
$ws.Range("M45").Value = -979
$ws.Range("H61").Value = 4153.46
$ws.Range("I61").Value = 2847.0476
$ws.Range("J61").Value = 5099.483
$ws.Range("K61").Value = 2847.0476
$ws.Range("L61").Value = 5099.483
$ws.Range("M61").Value = -2635.0476
$ws.Range("N61").Value = -5523.483
$ws.Range("H74").Value = 338985.7
$ws.Range("I74").Value = 437554.25
$ws.Range("J74").Value = 92564.3
$ws.Range("K74").Value = 437554.25
$ws.Range("L74").Value = 92564.3
$ws.Range("M74").Value = -436680.25
$ws.Range("N74").Value = -94312.3
$ws.Range("H77").Value = 338985.7
$ws.Range("I77").Value = 437554.25
$ws.Range("J77").Value = 92564.3
$ws.Range("K77").Value = 2187771.25
$ws.Range("L77").Value = 462821.5
$ws.Range("M77").Value = -2183403.25
$ws.Range("N77").Value = -471557.5
$ws.Range("H95").Value = 18569.334
$ws.Range("J95").Value = 18569.334
$ws.Range("L95").Value = 18569.334
$ws.Range("N95").Value = -24061.334
$ws.Range("H97").Value = 35715316
$ws.Range("I97").Value = 694.2273
$ws.Range("K97").Value = 694.2273
$ws.Range("M97").Value = -198.2273
$ws.Range("H110").Value = 260545.86
$ws.Range("I110").Value = 362764.2
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 362764.2
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -360719.2
$ws.Range("N110").Value = -9090
$ws.Range("H116").Value = 2016.4231
$ws.Range("I116").Value = 2004.3334
$ws.Range("J116").Value = 2043.625
$ws.Range("K116").Value = 2004.3334
$ws.Range("L116").Value = 2043.625
$ws.Range("M116").Value = 289.6666
$ws.Range("N116").Value = -6631.625
$ws.Range("H122").Value = 24232.732
$ws.Range("I122").Value = 2122.7
$ws.Range("J122").Value = 84532.82000000001
$ws.Range("K122").Value = 6368.099999999999
$ws.Range("L122").Value = 253598.46
$ws.Range("M122").Value = -3918.099999999999
$ws.Range("N122").Value = -258498.46
$ws.Range("H128").Value = 144999
$ws.Range("J128").Value = 144999
$ws.Range("L128").Value = 144999
$ws.Range("N128").Value = -154959
$ws.Range("H132").Value = 2231.5745
$ws.Range("I132").Value = 2196.7715
$ws.Range("J132").Value = 2333.0833
$ws.Range("K132").Value = 6590.314499999999
$ws.Range("L132").Value = 6999.249899999999
$ws.Range("M132").Value = -4060.314499999999
$ws.Range("N132").Value = -12059.2499
$ws.Range("H133").Value = 85500
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H136").Value = 4153.46
$ws.Range("I136").Value = 2847.0476
$ws.Range("J136").Value = 5099.483
$ws.Range("K136").Value = 8541.1428
$ws.Range("L136").Value = 15298.449
$ws.Range("M136").Value = -5991.1428
$ws.Range("N136").Value = -20398.449

# --- Sheet 3: BSM (59 cell updates) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 2016.4231
$ws.Range("I3").Value = 2004.3334
$ws.Range("J3").Value = 2043.625
$ws.Range("K3").Value = 2004.3334
$ws.Range("L3").Value = 2043.625
$ws.Range("M3").Value = -1890.3334
$ws.Range("N3").Value = -2271.625
$ws.Range("H20").Value = 1427.6428
$ws.Range("I20").Value = 1559.85
$ws.Range("J20").Value = 1097.125
$ws.Range("K20").Value = 1559.85
$ws.Range("L20").Value = 1097.125
$ws.Range("M20").Value = -1312.85
$ws.Range("N20").Value = -1591.125
$ws.Range("H75").Value = 24937.5
$ws.Range("J75").Value = 24875
$ws.Range("L75").Value = 24875
$ws.Range("N75").Value = -26747
$ws.Range("H76").Value = 78876.125
$ws.Range("J76").Value = 78876.125
$ws.Range("L76").Value = 78876.125
$ws.Range("N76").Value = -79506.125
$ws.Range("H78").Value = 24937.5
$ws.Range("J78").Value = 24875
$ws.Range("L78").Value = 74625
$ws.Range("N78").Value = -83985
$ws.Range("H79").Value = 78876.125
$ws.Range("J79").Value = 78876.125
$ws.Range("L79").Value = 78876.125
$ws.Range("N79").Value = -81060.125
$ws.Range("H94").Value = 65067.285
$ws.Range("I94").Value = 744.2
$ws.Range("J94").Value = 225875
$ws.Range("K94").Value = 744.2
$ws.Range("L94").Value = 225875
$ws.Range("M94").Value = -293.2
$ws.Range("N94").Value = -226777
$ws.Range("H105").Value = 4208.5557
$ws.Range("I105").Value = 3855.4
$ws.Range("J105").Value = 4650
$ws.Range("K105").Value = 3855.4
$ws.Range("L105").Value = 4650
$ws.Range("M105").Value = -2108.4
$ws.Range("N105").Value = -8144
$ws.Range("H107").Value = 2755.9333
$ws.Range("I107").Value = 3133.92
$ws.Range("J107").Value = 866
$ws.Range("K107").Value = 3133.92
$ws.Range("L107").Value = 866
$ws.Range("M107").Value = -1213.92
$ws.Range("N107").Value = -4706
$ws.Range("H134").Value = 2045.5536
$ws.Range("I134").Value = 1755.902
$ws.Range("K134").Value = 5267.706
$ws.Range("M134").Value = -2732.706
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# --- Sheet 4: CRP (87 cell updates) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 203.46153
$ws.Range("I7").Value = 115
$ws.Range("J7").Value = 345
$ws.Range("K7").Value = 115
$ws.Range("L7").Value = 345
$ws.Range("M7").Value = -2
$ws.Range("N7").Value = -571
$ws.Range("H16").Value = 6930.3105
$ws.Range("I16").Value = 5345.615
$ws.Range("J16").Value = 8217.875
$ws.Range("K16").Value = 5345.615
$ws.Range("L16").Value = 8217.875
$ws.Range("M16").Value = -5058.615
$ws.Range("N16").Value = -8791.875
$ws.Range("H31").Value = 2088.525
$ws.Range("I31").Value = 1792.4166
$ws.Range("J31").Value = 4753.5
$ws.Range("K31").Value = 1792.4166
$ws.Range("L31").Value = 4753.5
$ws.Range("M31").Value = -1497.4166
$ws.Range("N31").Value = -5343.5
$ws.Range("H34").Value = 2088.525
$ws.Range("I34").Value = 1792.4166
$ws.Range("J34").Value = 4753.5
$ws.Range("K34").Value = 1792.4166
$ws.Range("L34").Value = 4753.5
$ws.Range("M34").Value = -1590.4166
$ws.Range("N34").Value = -5157.5
$ws.Range("H52").Value = 86490.336
$ws.Range("I52").Value = 75000
$ws.Range("K52").Value = 75000
$ws.Range("M52").Value = -74706
$ws.Range("H58").Value = 2319.2222
$ws.Range("I58").Value = 1689.909
$ws.Range("K58").Value = 1689.909
$ws.Range("M58").Value = -1486.909
$ws.Range("H99").Value = 4668
$ws.Range("I99").Value = 4826.091
$ws.Range("J99").Value = 3798.5
$ws.Range("K99").Value = 4826.091
$ws.Range("L99").Value = 3798.5
$ws.Range("M99").Value = -3328.091
$ws.Range("N99").Value = -6794.5
$ws.Range("H105").Value = 2046.2858
$ws.Range("I105").Value = 1999.909
$ws.Range("K105").Value = 1999.909
$ws.Range("M105").Value = -252.9090000000001
$ws.Range("H107").Value = 30332466
$ws.Range("I107").Value = 50045776
$ws.Range("J107").Value = 4299.5386
$ws.Range("K107").Value = 50045776
$ws.Range("L107").Value = 4299.5386
$ws.Range("M107").Value = -50043856
$ws.Range("N107").Value = -8139.5386
$ws.Range("H113").Value = 6930.3105
$ws.Range("I113").Value = 5345.615
$ws.Range("J113").Value = 8217.875
$ws.Range("K113").Value = 5345.615
$ws.Range("L113").Value = 8217.875
$ws.Range("M113").Value = -3175.615
$ws.Range("N113").Value = -12557.875
$ws.Range("H126").Value = 4668
$ws.Range("I126").Value = 4826.091
$ws.Range("J126").Value = 3798.5
$ws.Range("K126").Value = 14478.273
$ws.Range("L126").Value = 11395.5
$ws.Range("M126").Value = -12008.273
$ws.Range("N126").Value = -16335.5
$ws.Range("H132").Value = 6823.3657
$ws.Range("I132").Value = 2354.5
$ws.Range("K132").Value = 7063.5
$ws.Range("M132").Value = -4533.5
$ws.Range("H134").Value = 1084.7273
$ws.Range("I134").Value = 1063.4419
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 3190.3257
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -655.3257000000003
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 2319.2222
$ws.Range("I136").Value = 1689.909
$ws.Range("K136").Value = 5069.727000000001
$ws.Range("M136").Value = -2519.727000000001
$ws.Range("H138").Value = 68905.17999999999
$ws.Range("J138").Value = 68695.7
$ws.Range("L138").Value = 68695.7
$ws.Range("N138").Value = -78975.7

# --- Sheet 5: CUL (80 cell updates) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 49442428
$ws.Range("J4").Value = 128294620
$ws.Range("L4").Value = 384883860
$ws.Range("N4").Value = -384884084
$ws.Range("H12").Value = 447.9
$ws.Range("I12").Value = 296.4
$ws.Range("J12").Value = 599.4
$ws.Range("K12").Value = 889.1999999999999
$ws.Range("L12").Value = 1798.2
$ws.Range("M12").Value = -716.1999999999999
$ws.Range("N12").Value = -2144.2
$ws.Range("H26").Value = 1775.3636
$ws.Range("I26").Value = 2716.6667
$ws.Range("J26").Value = 645.8
$ws.Range("K26").Value = 8150.000100000001
$ws.Range("L26").Value = 1937.4
$ws.Range("M26").Value = -7862.000100000001
$ws.Range("N26").Value = -2513.4
$ws.Range("H34").Value = 2599.8
$ws.Range("I34").Value = 633.3333
$ws.Range("K34").Value = 1899.9999
$ws.Range("M34").Value = -1815.9999
$ws.Range("H62").Value = 9499
$ws.Range("H65").Value = 9499
$ws.Range("H80").Value = 53799.8
$ws.Range("I80").Value = 47999.668
$ws.Range("J80").Value = 62500
$ws.Range("K80").Value = 143999.004
$ws.Range("L80").Value = 187500
$ws.Range("M80").Value = -143063.004
$ws.Range("N80").Value = -189372
$ws.Range("H83").Value = 53799.8
$ws.Range("I83").Value = 47999.668
$ws.Range("J83").Value = 62500
$ws.Range("K83").Value = 431997.012
$ws.Range("L83").Value = 562500
$ws.Range("M83").Value = -427317.012
$ws.Range("N83").Value = -571860
$ws.Range("H92").Value = 358
$ws.Range("I92").Value = 99
$ws.Range("K92").Value = 297
$ws.Range("M92").Value = 951
$ws.Range("H113").Value = 1986
$ws.Range("I113").Value = 728.8
$ws.Range("K113").Value = 2186.4
$ws.Range("M113").Value = -16.39999999999964
$ws.Range("H121").Value = 2895786
$ws.Range("J121").Value = 3860927.5
$ws.Range("L121").Value = 11582782.5
$ws.Range("N121").Value = -11585402.5
$ws.Range("H122").Value = 696.0714
$ws.Range("I122").Value = 720.8
$ws.Range("J122").Value = 682.3333
$ws.Range("K122").Value = 6487.2
$ws.Range("L122").Value = 6140.9997
$ws.Range("M122").Value = -4037.2
$ws.Range("N122").Value = -11040.9997
$ws.Range("H129").Value = 2462.4
$ws.Range("I129").Value = 849.375
$ws.Range("J129").Value = 4305.857
$ws.Range("K129").Value = 2548.125
$ws.Range("L129").Value = 12917.571
$ws.Range("M129").Value = 2451.875
$ws.Range("N129").Value = -22917.571
$ws.Range("H131").Value = 5023.5454
$ws.Range("J131").Value = 5608.857
$ws.Range("L131").Value = 16826.571
$ws.Range("N131").Value = -26906.571
$ws.Range("H134").Value = 7137.75
$ws.Range("I134").Value = 7137.75
$ws.Range("K134").Value = 21413.25
$ws.Range("M134").Value = -16343.25
$ws.Range("H139").Value = 4533.7646
$ws.Range("I139").Value = 2710.4167
$ws.Range("K139").Value = 8131.250100000001
$ws.Range("M139").Value = -2991.250100000001
$ws.Range("H140").Value = 2194.4285
$ws.Range("I140").Value = 1473.6
$ws.Range("K140").Value = 4420.799999999999
$ws.Range("M140").Value = 759.2000000000007

# --- Sheet 6: GSM (68 cell updates) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H57").Value = 24973.666
$ws.Range("I57").Value = 18000
$ws.Range("J57").Value = 28460.5
$ws.Range("K57").Value = 18000
$ws.Range("L57").Value = 28460.5
$ws.Range("M57").Value = -17180
$ws.Range("N57").Value = -30100.5
$ws.Range("H70").Value = 7885.4546
$ws.Range("I70").Value = 6739.25
$ws.Range("J70").Value = 8540.429
$ws.Range("K70").Value = 6739.25
$ws.Range("L70").Value = 8540.429
$ws.Range("M70").Value = -6469.25
$ws.Range("N70").Value = -9080.429
$ws.Range("H73").Value = 7885.4546
$ws.Range("I73").Value = 6739.25
$ws.Range("J73").Value = 8540.429
$ws.Range("K73").Value = 6739.25
$ws.Range("L73").Value = 8540.429
$ws.Range("M73").Value = -5803.25
$ws.Range("N73").Value = -10412.429
$ws.Range("H93").Value = 59666.668
$ws.Range("J93").Value = 59666.668
$ws.Range("L93").Value = 59666.668
$ws.Range("N93").Value = -63410.668
$ws.Range("H97").Value = 43520100
$ws.Range("I97").Value = 58878976
$ws.Range("J97").Value = 3283.5
$ws.Range("K97").Value = 58878976
$ws.Range("L97").Value = 3283.5
$ws.Range("M97").Value = -58878480
$ws.Range("N97").Value = -4275.5
$ws.Range("H102").Value = 26739.61
$ws.Range("I102").Value = 33826.25
$ws.Range("J102").Value = 10541.571
$ws.Range("K102").Value = 33826.25
$ws.Range("L102").Value = 10541.571
$ws.Range("M102").Value = -32204.25
$ws.Range("N102").Value = -13785.571
$ws.Range("H113").Value = 2318.1538
$ws.Range("I113").Value = 2343
$ws.Range("J113").Value = 2181.5
$ws.Range("K113").Value = 2343
$ws.Range("L113").Value = 2181.5
$ws.Range("M113").Value = -173
$ws.Range("N113").Value = -6521.5
$ws.Range("H122").Value = 60495.234
$ws.Range("I122").Value = 72979.92999999999
$ws.Range("J122").Value = 2233.3333
$ws.Range("K122").Value = 218939.79
$ws.Range("L122").Value = 6699.999899999999
$ws.Range("M122").Value = -216489.79
$ws.Range("N122").Value = -11599.9999
$ws.Range("H126").Value = 55358.445
$ws.Range("I126").Value = 61539.062
$ws.Range("K126").Value = 184617.186
$ws.Range("M126").Value = -182147.186
$ws.Range("H132").Value = 3202.7778
$ws.Range("I132").Value = 3789.2104
$ws.Range("J132").Value = 1810
$ws.Range("K132").Value = 11367.6312
$ws.Range("L132").Value = 5430
$ws.Range("M132").Value = -8837.6312
$ws.Range("N132").Value = -10490
$ws.Range("H136").Value = 20105.572
$ws.Range("J136").Value = 20105.572
$ws.Range("L136").Value = 60316.716
$ws.Range("N136").Value = -65416.716

# --- Sheet 7: LTW (111 cell updates) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 21092.785
$ws.Range("I7").Value = 30585.389
$ws.Range("J7").Value = 4006.1
$ws.Range("K7").Value = 30585.389
$ws.Range("L7").Value = 4006.1
$ws.Range("M7").Value = -30473.389
$ws.Range("N7").Value = -4230.1
$ws.Range("H16").Value = 1235.2245
$ws.Range("I16").Value = 1203.7368
$ws.Range("K16").Value = 1203.7368
$ws.Range("M16").Value = -1033.7368
$ws.Range("H22").Value = 4258.6665
$ws.Range("I22").Value = 3016.8
$ws.Range("J22").Value = 4879.6
$ws.Range("K22").Value = 3016.8
$ws.Range("L22").Value = 4879.6
$ws.Range("M22").Value = -2721.8
$ws.Range("N22").Value = -5469.6
$ws.Range("H26").Value = 11022.8
$ws.Range("I26").Value = 11022.8
$ws.Range("K26").Value = 11022.8
$ws.Range("M26").Value = -10727.8
$ws.Range("H27").Value = 4258.6665
$ws.Range("I27").Value = 3016.8
$ws.Range("J27").Value = 4879.6
$ws.Range("K27").Value = 3016.8
$ws.Range("L27").Value = 4879.6
$ws.Range("M27").Value = -2909.8
$ws.Range("N27").Value = -5093.6
$ws.Range("H40").Value = 28470.824
$ws.Range("I40").Value = 35362
$ws.Range("K40").Value = 35362
$ws.Range("M40").Value = -35226
$ws.Range("H46").Value = 2885.9692
$ws.Range("J46").Value = 2979.6453
$ws.Range("L46").Value = 2979.6453
$ws.Range("N46").Value = -3355.6453
$ws.Range("H68").Value = 3149.7778
$ws.Range("I68").Value = 1192.8572
$ws.Range("J68").Value = 9999
$ws.Range("K68").Value = 1192.8572
$ws.Range("L68").Value = 9999
$ws.Range("M68").Value = -443.8571999999999
$ws.Range("N68").Value = -11497
$ws.Range("H71").Value = 3149.7778
$ws.Range("I71").Value = 1192.8572
$ws.Range("J71").Value = 9999
$ws.Range("K71").Value = 5964.286
$ws.Range("L71").Value = 49995
$ws.Range("M71").Value = -2220.286
$ws.Range("N71").Value = -57483
$ws.Range("H82").Value = 1176.7
$ws.Range("I82").Value = 1047.5
$ws.Range("J82").Value = 1693.5
$ws.Range("K82").Value = 1047.5
$ws.Range("L82").Value = 1693.5
$ws.Range("M82").Value = -686.5
$ws.Range("N82").Value = -2415.5
$ws.Range("H85").Value = 1176.7
$ws.Range("I85").Value = 1047.5
$ws.Range("J85").Value = 1693.5
$ws.Range("K85").Value = 1047.5
$ws.Range("L85").Value = 1693.5
$ws.Range("M85").Value = 200.5
$ws.Range("N85").Value = -4189.5
$ws.Range("H93").Value = 2024.2
$ws.Range("I93").Value = 2255.25
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 2255.25
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = -1007.25
$ws.Range("N93").Value = -3596
$ws.Range("H100").Value = 2360.5557
$ws.Range("I100").Value = 2097.1875
$ws.Range("J100").Value = 4467.5
$ws.Range("K100").Value = 2097.1875
$ws.Range("L100").Value = 4467.5
$ws.Range("M100").Value = -1556.1875
$ws.Range("N100").Value = -5549.5
$ws.Range("H122").Value = 3600
$ws.Range("I122").Value = 3575.15
$ws.Range("J122").Value = 3765.6667
$ws.Range("K122").Value = 10725.45
$ws.Range("L122").Value = 11297.0001
$ws.Range("M122").Value = -8275.450000000001
$ws.Range("N122").Value = -16197.0001
$ws.Range("H126").Value = 21092.785
$ws.Range("I126").Value = 30585.389
$ws.Range("J126").Value = 4006.1
$ws.Range("K126").Value = 91756.167
$ws.Range("L126").Value = 12018.3
$ws.Range("M126").Value = -89286.167
$ws.Range("N126").Value = -16958.3
$ws.Range("H128").Value = 84714
$ws.Range("J128").Value = 84714
$ws.Range("L128").Value = 84714
$ws.Range("N128").Value = -94674
$ws.Range("H132").Value = 2632.2827
$ws.Range("I132").Value = 1993.6364
$ws.Range("J132").Value = 4253.4614
$ws.Range("K132").Value = 5980.9092
$ws.Range("L132").Value = 12760.3842
$ws.Range("M132").Value = -3450.9092
$ws.Range("N132").Value = -17820.3842
$ws.Range("H136").Value = 22199.238
$ws.Range("I136").Value = 2289.2778
$ws.Range("J136").Value = 93875.10000000001
$ws.Range("K136").Value = 6867.8334
$ws.Range("L136").Value = 281625.3
$ws.Range("M136").Value = -4317.8334
$ws.Range("N136").Value = -286725.3

# --- Sheet 8: WVR (62 cell updates) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H14").Value = 4897.6
$ws.Range("I14").Value = 498
$ws.Range("J14").Value = 11497
$ws.Range("K14").Value = 498
$ws.Range("L14").Value = 11497
$ws.Range("M14").Value = -330
$ws.Range("N14").Value = -11833
$ws.Range("H81").Value = 1331.375
$ws.Range("I81").Value = 1331.375
$ws.Range("K81").Value = 2662.75
$ws.Range("M81").Value = -1601.75
$ws.Range("H84").Value = 1331.375
$ws.Range("I84").Value = 1331.375
$ws.Range("K84").Value = 13313.75
$ws.Range("M84").Value = -8009.75
$ws.Range("H96").Value = 6113.8335
$ws.Range("I96").Value = 8330.666999999999
$ws.Range("J96").Value = 3897
$ws.Range("K96").Value = 8330.666999999999
$ws.Range("L96").Value = 3897
$ws.Range("M96").Value = -6957.666999999999
$ws.Range("N96").Value = -6643
$ws.Range("H102").Value = 149998.5
$ws.Range("J102").Value = 149998.5
$ws.Range("L102").Value = 149998.5
$ws.Range("N102").Value = -156488.5
$ws.Range("H106").Value = 48266.332
$ws.Range("J106").Value = 48266.332
$ws.Range("L106").Value = 48266.332
$ws.Range("N106").Value = -50790.332
$ws.Range("H113").Value = 1287.04
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 910.2222
$ws.Range("K113").Value = 4497
$ws.Range("L113").Value = 2730.6666
$ws.Range("M113").Value = -2327
$ws.Range("N113").Value = -7070.6666
$ws.Range("H122").Value = 3249.0688
$ws.Range("I122").Value = 3076.652
$ws.Range("K122").Value = 9229.956
$ws.Range("M122").Value = -6779.956
$ws.Range("H126").Value = 71436320
$ws.Range("I126").Value = 111118584
$ws.Range("J126").Value = 8238.799999999999
$ws.Range("K126").Value = 333355752
$ws.Range("L126").Value = 24716.4
$ws.Range("M126").Value = -333353282
$ws.Range("N126").Value = -29656.4
$ws.Range("H132").Value = 2596.551
$ws.Range("I132").Value = 2541.238
$ws.Range("J132").Value = 2928.4285
$ws.Range("K132").Value = 7623.714
$ws.Range("L132").Value = 8785.2855
$ws.Range("M132").Value = -5093.714
$ws.Range("N132").Value = -13845.2855
$ws.Range("H136").Value = 58210.117
$ws.Range("I136").Value = 38544.23
$ws.Range("J136").Value = 122124.25
$ws.Range("K136").Value = 115632.69
$ws.Range("L136").Value = 366372.75
$ws.Range("M136").Value = -113082.69
$ws.Range("N136").Value = -371472.75
